# Apply the commit "drag and drop actions":
#  1. Re-cache the datetimeFigureOut date placeholders (slide master + every
#     slide layout) from 8/4/13 to 8/11/13.
#  2. Resize/reposition/restyle the "Hello World" textbox on slide 1 -
#     move+grow it, give it a yellow fill, and bump the run font size.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Date placeholders on the slide master and every slide layout.
# ---------------------------------------------------------------------
function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "8/11/13"
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-DatePlaceholders $layout.Shapes
}

# ---------------------------------------------------------------------
# 2. "Hello World" textbox on slide 1.
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)

# Position / size (EMU -> points, 12700 EMU per point). A tiny epsilon is
# added before conversion so the host's point->EMU rounding lands on the
# exact target EMU instead of one unit short.
$emuPerPt = 12700
$eps = 0.00005

$sh.Left   = (2357300 / $emuPerPt) + $eps
$sh.Top    = (2654478 / $emuPerPt) + $eps
$sh.Width  = (4292611 / $emuPerPt) + $eps
$sh.Height = (1107996 / $emuPerPt) + $eps

# Fill: noFill -> solid yellow (FFFF00).
$sh.Fill.Solid()
$sh.Fill.ForeColor.RGB = 65535

# Font size: 36pt -> 66pt.
$sh.TextFrame.TextRange.Font.Size = 66
